# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (totals) sheet,
#    and populate it with the same per-fund column layout used by the other
#    quarterly sheets (2020-Q4 / 2021-Q3 / 2021-Q4).
# 2. Insert a new leading data row in "总计" summarising the new quarter and
#    renumber the existing index column (A) beneath it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet, placed immediately before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Match the look of the other quarterly sheets: bold/centered/bordered
# header row and bold/centered index column (style carried via the
# existing "2021-Q4" sheet rather than rebuilt by hand).
$prior = $wb.Worksheets.Item("2021-Q4")
$prior.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$prior.Range("A2:A3").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# These columns hold text-formatted numbers in the source data (fund code
# needs its leading zero preserved; the scale/position figures are stored
# as plain text too), so force text formatting before writing the values.
$q1.Range("B2:B3").NumberFormat = "@"
$q1.Range("D2:G3").NumberFormat = "@"

# Row 2 - 006792
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "006792"
$q1.Range("C2").Value = "鹏华香港美国互联网股票（LOF）美元现汇"
$q1.Range("D2").Value = "1.43"
$q1.Range("E2").Value = "83.72"
$q1.Range("F2").Value = "6.79"
$q1.Range("G2").Value = "0.0971"
$q1.Range("H2").Value = 3

# Row 3 - 160644
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "160644"
$q1.Range("C3").Value = "鹏华香港美国互联网股票（LOF）人民币"
$q1.Range("D3").Value = "1.43"
$q1.Range("E3").Value = "83.72"
$q1.Range("F3").Value = "6.79"
$q1.Range("G3").Value = "0.0971"
$q1.Range("H3").Value = 3

# The "@" number format served its purpose (locking the values in as text);
# drop it again so the cells end up unstyled, matching the plain data cells
# on the sibling quarterly sheets. Their stored type stays Text.
$q1.Range("B2:B3").ClearFormats()
$q1.Range("D2:G3").ClearFormats()

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet with the new quarter on top
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# The row-insert bleeds the header's bold/border formatting into B2:D2 and
# drops the index-column style from A2 - restore both to match the rest of
# column A / the plain data cells below.
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.19

# Renumber the index column for the rows pushed down
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
